$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Date" column for this export was populated with the wrong value
# ("5-23-2011-12" - a mangled form of the source NBA-stats date) for every
# data row. Correct it to the real ISO date, "2012-05-23".
$oldDate = "5-23-2011-12"
$newDate = "2012-05-23"

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$firstCol = $usedRange.Column
$lastRow = $firstRow + $usedRange.Rows.Count - 1
$lastCol = $firstCol + $usedRange.Columns.Count - 1

# Locate the "Date" header column dynamically (it's BF / column 58 in this
# sheet) rather than hard-coding it.
$dateCol = -1
for ($c = $firstCol; $c -le $lastCol; $c++) {
    if ($ws.Cells.Item($firstRow, $c).Value2 -eq "Date") {
        $dateCol = $c
        break
    }
}

if ($dateCol -ne -1) {
    for ($row = $firstRow + 1; $row -le $lastRow; $row++) {
        $cell = $ws.Cells.Item($row, $dateCol)
        if ($cell.Value2 -eq $oldDate) {
            # A plain string that looks like a date (e.g. "2012-05-23") gets
            # silently auto-converted to a date serial number when assigned
            # to .Value, so force the cell to text first, write the literal
            # string, then clear the formatting residue back to General so
            # no stray number-format style is left on the cell.
            $cell.NumberFormat = "@"
            $cell.Value = $newDate
            $cell.ClearFormats()
        }
    }
}
